{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// The paragraph to remove is the one asking about the sequence diagram\n// (\"Diagramma di sequenza ...\") \u2014 match on a stable ASCII-only prefix so\n// we don't depend on exact curly-quote / whitespace rendering.\nconst marker = \"Diagramma di sequenza\";\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(marker) !== -1) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the question paragraph about the sequence diagram\n# (\"Diagramma di sequenza ...\"). Match on a stable ASCII-only substring\n# so curly-quote/whitespace rendering differences don't matter.\nforeach ($p in @($d.Paragraphs)) {\n    if ($p.Range.Text -like \"*Diagramma di sequenza*\") {\n        $p.Range.Delete()\n    }\n}\n"}
